# Equipment Manage: rename the "Sucursales" worksheet to "Equipos".
# Renaming the sheet also updates the "sucursales" defined name, whose
# formula references the sheet by name (Sucursales!$A$4:$H$5 -> Equipos!$A$4:$H$5).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Equipos"
